# Apply the "Fixed update to excel issue" changes:
#  1. Rename "Requested quantity" header -> "Weekly_PO_Qty" on "Weekly Quantity" sheet
#  2. Rename "Requested quantity" header -> "Monthly_PO_Qty" on "Monthly Trend" sheet
#  3. Add a new "PO Forecast" sheet (after "Monthly Trend") with forecast data

$wb = $excel.ActiveWorkbook

# --- 1 & 2: rename the "Requested quantity" headers -----------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: add the new "PO Forecast" sheet as the last tab --------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the bold/centered/bordered header style already used elsewhere in
# the workbook (copy format only, so the existing style is reused instead
# of a new one being minted).
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Data rows
$data = @(
  @(45053.99999999999, 1,  1.000003827037351,  1.000003827063771),
  @(45067.99999999999, 3,  3.000004824260171,  3.00000482428732),
  @(45074.99999999999, 4,  4.000005315287569,  4.000005330275259),
  @(45081.99999999999, 5,  5.00000579716235,   5.000005847144346),
  @(45088.99999999999, 6,  6.000006271876677,  6.000006369493723),
  @(45095.99999999999, 7,  7.000006741935914,  7.00000689450148),
  @(45102.99999999999, 8,  8.000007211234664,  8.000007425811182),
  @(45109.99999999999, 9,  9.000007676001367,  9.00000795802088),
  @(45116.99999999999, 10, 10.00000813169523,  10.00000849032945),
  @(45123.99999999999, 11, 11.00000858357773,  11.00000903707643)
)

$row = 2
foreach ($d in $data) {
    $wsForecast.Range("A$row").Value = $d[0]
    $wsForecast.Range("B$row").Value = $d[1]
    $wsForecast.Range("C$row").Value = $d[2]
    $wsForecast.Range("D$row").Value = $d[3]
    $row++
}

# Match the date-time number format used for the "ds" column elsewhere.
$wsForecast.Range("A2:A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$wsForecast.Range("A1").Select()
